$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5744
$ws.Range("K3").Value = 8181
$ws.Range("L3").Value = 6266
$ws.Range("L4").Value = 1545
$ws.Range("L5").Value = 374
$ws.Range("L6").Value = 5143
$ws.Range("K7").Value = 27585
$ws.Range("L7").Value = 19072

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 61
$ws.Range("L3").Value = 53
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 207

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 447
$ws.Range("L7").Value = 1258

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 236
$ws.Range("L6").Value = 248
$ws.Range("L7").Value = 868

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 217
$ws.Range("L3").Value = 256
$ws.Range("L6").Value = 197
$ws.Range("L7").Value = 729

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 121
$ws.Range("L7").Value = 367

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 94
$ws.Range("L7").Value = 329

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L8").Value = 1258
$ws.Range("L10").Value = 128
$ws.Range("L11").Value = 314
$ws.Range("L18").Value = 132
$ws.Range("L19").Value = 521
$ws.Range("L20").Value = 474
$ws.Range("L27").Value = 169
$ws.Range("K29").Value = 1520
$ws.Range("L29").Value = 1077
$ws.Range("L33").Value = 868
$ws.Range("L37").Value = 729
$ws.Range("L39").Value = 11
$ws.Range("L42").Value = 614
$ws.Range("L47").Value = 128
$ws.Range("L48").Value = 250
$ws.Range("L50").Value = 94
$ws.Range("L51").Value = 239
$ws.Range("L52").Value = 396
$ws.Range("L53").Value = 207
$ws.Range("K59").Value = 44
$ws.Range("L59").Value = 33
$ws.Range("K63").Value = 177
$ws.Range("L65").Value = 367
$ws.Range("L67").Value = 665
$ws.Range("L70").Value = 34
$ws.Range("L77").Value = 129
$ws.Range("L79").Value = 524
$ws.Range("L85").Value = 945
$ws.Range("L86").Value = 128
$ws.Range("L88").Value = 203
$ws.Range("L91").Value = 257
$ws.Range("L96").Value = 218
$ws.Range("L97").Value = 155
$ws.Range("L98").Value = 102
$ws.Range("L99").Value = 329
$ws.Range("K101").Value = 27585
$ws.Range("L101").Value = 19072

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L4").Value = 45
$ws.Range("L5").Value = 18
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 665

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 319
$ws.Range("L3").Value = 413
$ws.Range("K4").Value = 74
$ws.Range("L6").Value = 268
$ws.Range("K7").Value = 1520
$ws.Range("L7").Value = 1077

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L4").Value = 50
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 250

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 186
$ws.Range("L7").Value = 521

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 168
$ws.Range("L7").Value = 614

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 49
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 119
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 257

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 167
$ws.Range("L6").Value = 140
$ws.Range("L7").Value = 524

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 162
$ws.Range("L7").Value = 474

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 50
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L3").Value = 3
$ws.Range("L6").Value = 11

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 118
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 314

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K3").Value = 13
$ws.Range("L4").Value = 3
$ws.Range("K7").Value = 44
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L4").Value = 13
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L3").Value = 26
$ws.Range("L4").Value = 68
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 239

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 286
$ws.Range("L3").Value = 390
$ws.Range("L6").Value = 195
$ws.Range("L7").Value = 945

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 125
$ws.Range("L3").Value = 127
$ws.Range("L5").Value = 11
$ws.Range("L7").Value = 396
